# Generate Report for Handoff
# Adds two new rows (for d7cf9305-...-fb7da0.md and dcd52aba-...-fe422d.md)
# to the Overview / zh-cn / de-de sheets, growing every table from 3 -> 5 rows.

$wb = $excel.ActiveWorkbook

$HYPER_BLUE = 15570276   # BGR encoding of RGB FF6495ED (the workbook's custom hyperlink color)

function Style-AsLink($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = 2
    $rng.Font.Color = $HYPER_BLUE
}

# ---------------------------------------------------------------------------
# Sheet "Overview" (sheet1) - columns A..G
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4,1).Value = "d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md", "", "", "e2e\d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md")
Style-AsLink $wsOverview.Range("B4")
$wsOverview.Cells.Item(4,3).Value = ".md"
$wsOverview.Cells.Item(4,5).Value = "Ready for handoff"
$wsOverview.Cells.Item(4,6).Value = "Ready for handoff"
$wsOverview.Cells.Item(4,7).Value = "2016-08-27 12:38:17"
$wsOverview.Range("G4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$wsOverview.Cells.Item(5,1).Value = "dcd52aba-f709-418b-91ef-74ba14fe422d.md"
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/dcd52aba-f709-418b-91ef-74ba14fe422d.md", "", "", "e2e\dcd52aba-f709-418b-91ef-74ba14fe422d.md")
Style-AsLink $wsOverview.Range("B5")
$wsOverview.Cells.Item(5,3).Value = ".md"
$wsOverview.Cells.Item(5,5).Value = "Ready for handoff"
$wsOverview.Cells.Item(5,6).Value = "Ready for handoff"
$wsOverview.Cells.Item(5,7).Value = "2016-08-27 12:38:17"
$wsOverview.Range("G5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn" (sheet2) - columns A..P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md", "", "", "d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md")
Style-AsLink $wsZh.Range("A4")
$wsZh.Cells.Item(4,2).Value = ".md"
$wsZh.Cells.Item(4,3).Value = "Ready for handoff"
$wsZh.Cells.Item(4,4).Value = "e2e"
$wsZh.Cells.Item(4,5).Value = "ht"
$wsZh.Cells.Item(4,6).Value = "False"
$wsZh.Cells.Item(4,7).Value = "d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.7cdd7cee52d46e5d0524ed399816b3926e64531c.zh-cn.xlf"
$wsZh.Cells.Item(4,8).Value = "2016-08-27 12:38:11"
$wsZh.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$wsZh.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,13).Value = "True"
$wsZh.Cells.Item(4,15).Value = "False"

$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/dcd52aba-f709-418b-91ef-74ba14fe422d.md", "", "", "dcd52aba-f709-418b-91ef-74ba14fe422d.md")
Style-AsLink $wsZh.Range("A5")
$wsZh.Cells.Item(5,2).Value = ".md"
$wsZh.Cells.Item(5,3).Value = "Ready for handoff"
$wsZh.Cells.Item(5,4).Value = "e2e"
$wsZh.Cells.Item(5,5).Value = "ht"
$wsZh.Cells.Item(5,6).Value = "False"
$wsZh.Cells.Item(5,7).Value = "dcd52aba-f709-418b-91ef-74ba14fe422d.2c90cd0b9a6bf19fee2b3e2d4913e4eb0aab0d7b.zh-cn.xlf"
$wsZh.Cells.Item(5,8).Value = "2016-08-27 12:38:11"
$wsZh.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$wsZh.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(5,13).Value = "True"
$wsZh.Cells.Item(5,15).Value = "False"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet "de-de" (sheet3) - columns A..P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md", "", "", "d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md")
Style-AsLink $wsDe.Range("A4")
$wsDe.Cells.Item(4,2).Value = ".md"
$wsDe.Cells.Item(4,3).Value = "Ready for handoff"
$wsDe.Cells.Item(4,4).Value = "e2e"
$wsDe.Cells.Item(4,5).Value = "ht"
$wsDe.Cells.Item(4,6).Value = "False"
$wsDe.Cells.Item(4,7).Value = "d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.7cdd7cee52d46e5d0524ed399816b3926e64531c.de-de.xlf"
$wsDe.Cells.Item(4,8).Value = "2016-08-27 12:38:17"
$wsDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$wsDe.Range("K4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,13).Value = "True"
$wsDe.Cells.Item(4,15).Value = "False"

$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/dcd52aba-f709-418b-91ef-74ba14fe422d.md", "", "", "dcd52aba-f709-418b-91ef-74ba14fe422d.md")
Style-AsLink $wsDe.Range("A5")
$wsDe.Cells.Item(5,2).Value = ".md"
$wsDe.Cells.Item(5,3).Value = "Ready for handoff"
$wsDe.Cells.Item(5,4).Value = "e2e"
$wsDe.Cells.Item(5,5).Value = "ht"
$wsDe.Cells.Item(5,6).Value = "False"
$wsDe.Cells.Item(5,7).Value = "dcd52aba-f709-418b-91ef-74ba14fe422d.2c90cd0b9a6bf19fee2b3e2d4913e4eb0aab0d7b.de-de.xlf"
$wsDe.Cells.Item(5,8).Value = "2016-08-27 12:38:17"
$wsDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(5,11).Value = "0001-01-01 00:00:00"
$wsDe.Range("K5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(5,13).Value = "True"
$wsDe.Cells.Item(5,15).Value = "False"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P5"))

Write-Host "Report rows added for d7cf9305-6c7b-4efe-a2a9-3ef870fb7da0.md and dcd52aba-f709-418b-91ef-74ba14fe422d.md"
